$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "https://www.bidtheatre.com/"
$ws.Range("A5").Value = "chrome-extension://cjpalhdlnbpafiamejdnhcphjbkeiagm/document-blocked.html?details=%7B%22url%22%3A%22https%3A%2F%2Fwww.bidtheatre.com%2F%22%2C%22hn%22%3A%22www.bidtheatre.com%22%2C%22dn%22%3A%22bidtheatre.com%22%2C%22fs%22%3A%22%7C%7Cbidtheatre.com%5E%22%7D"
